$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "60.768.47"
$ws.Range("E2").Value = "  -0.19%  "

# Row 3
$ws.Range("D3").Value = "2.343.75"
$ws.Range("E3").Value = "  -1.80%  "

# Row 5
$ws.Range("D5").Formula = '="543.43"'
$ws.Range("D5").Copy() | Out-Null
$ws.Range("D5").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("E5").Value = "  -0.12%  "

# Row 6
$ws.Range("D6").Formula = '="136.35"'
$ws.Range("D6").Copy() | Out-Null
$ws.Range("D6").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("E6").Value = "  -3.26%  "

# Row 7
$ws.Range("E7").Value = "  +0.00%  "

# Row 8
$ws.Range("D8").Formula = '="0.525"'
$ws.Range("D8").Copy() | Out-Null
$ws.Range("D8").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("E8").Value = "  -9.20%  "

# Row 9
$ws.Range("D9").Value = "2.342.43"
$ws.Range("E9").Value = "  -1.72%  "

# Row 10
$ws.Range("E10").Value = "  -0.77%  "

# Row 11
$ws.Range("E11").Value = "  +1.66%  "

# Row 12
$ws.Range("E12").Value = "  -1.36%  "

# Row 13
$ws.Range("E13").Value = "  -0.57%  "

# Row 14
$ws.Range("E14").Value = "  -3.50%  "

# Row 15
$ws.Range("D15").Value = "2.765.61"
$ws.Range("E15").Value = "  -1.76%  "

# Row 16
$ws.Range("D16").Value = "60.735.24"
$ws.Range("E16").Value = "  +0.28%  "

# Row 17
$ws.Range("E17").Value = "  -2.97%  "

# Row 18
$ws.Range("D18").Value = "2.342.34"
$ws.Range("E18").Value = "  -1.73%  "

# Row 19
$ws.Range("D19").Formula = '="10.61"'
$ws.Range("D19").Copy() | Out-Null
$ws.Range("D19").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("E19").Value = "  -0.17%  "

# Row 20
$ws.Range("E20").Value = "  +0.01%  "

# Row 21
$ws.Range("D21").Formula = '="318.46"'
$ws.Range("D21").Copy() | Out-Null
$ws.Range("D21").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("E21").Value = "  +0.12%  "

# Row 22
$ws.Range("E22").Value = "  -2.84%  "

# Row 23
$ws.Range("E23").Value = "  +0.00%  "

# Row 24
$ws.Range("D24").Formula = '="63.23"'
$ws.Range("D24").Copy() | Out-Null
$ws.Range("D24").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("E24").Value = "  +0.04%  "

# Row 25
$ws.Range("D25").Formula = '="1.70"'
$ws.Range("D25").Copy() | Out-Null
$ws.Range("D25").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("E25").Value = "  -6.56%  "

# Row 26
$ws.Range("D26").Formula = '="8.48"'
$ws.Range("D26").Copy() | Out-Null
$ws.Range("D26").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("E26").Value = "  +8.32%  "

# Row 27
$ws.Range("E27").Value = "  +0.27%  "

# Row 28
$ws.Range("D28").Value = "2.458.02"
$ws.Range("E28").Value = "  -1.64%  "

# Row 29
$ws.Range("D29").Formula = '="7.92"'
$ws.Range("D29").Copy() | Out-Null
$ws.Range("D29").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("E29").Value = "  -1.26%  "

# Row 30
$ws.Range("B30").Value = "Bittensor"
$ws.Range("C30").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D30").Formula = '="497.54"'
$ws.Range("D30").Copy() | Out-Null
$ws.Range("D30").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("E30").Value = "  -5.01%  "

# Row 31
$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").Formula = '="1.37"'
$ws.Range("D31").Copy() | Out-Null
$ws.Range("D31").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("E31").Value = "  -4.37%  "

# Row 32
$ws.Range("E32").Value = "  -8.06%  "

# Row 33
$ws.Range("E33").Value = "  +0.51%  "

# Row 34
$ws.Range("D34").Formula = '="1.78"'
$ws.Range("D34").Copy() | Out-Null
$ws.Range("D34").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("E34").Value = "  -2.33%  "

# Row 35
$ws.Range("E35").Value = "  -4.35%  "

# Row 36
$ws.Range("E36").Value = "  +0.05%  "

# Row 37
$ws.Range("E37").Value = "  -1.93%  "

# Row 38
$ws.Range("E38").Value = "  -0.47%  "

# Row 39
$ws.Range("D39").Formula = '="18.47"'
$ws.Range("D39").Copy() | Out-Null
$ws.Range("D39").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("E39").Value = "  +1.95%  "

# Row 40
$ws.Range("D40").Formula = '="1.82"'
$ws.Range("D40").Copy() | Out-Null
$ws.Range("D40").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("E40").Value = "  +5.05%  "

# Row 41
$ws.Range("D41").Formula = '="5.24"'
$ws.Range("D41").Copy() | Out-Null
$ws.Range("D41").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("E41").Value = "  -4.58%  "

# Row 42
$ws.Range("D42").Formula = '="142.75"'
$ws.Range("D42").Copy() | Out-Null
$ws.Range("D42").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("E42").Value = "  +3.17%  "

# Row 43
$ws.Range("E43").Value = "  -0.09%  "

# Row 44
$ws.Range("D44").Formula = '="40.51"'
$ws.Range("D44").Copy() | Out-Null
$ws.Range("D44").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("E44").Value = "  +0.84%  "

# Row 45
$ws.Range("D45").Formula = '="142.37"'
$ws.Range("D45").Copy() | Out-Null
$ws.Range("D45").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("E45").Value = "  +1.44%  "

# Row 46
$ws.Range("D46").Formula = '="3.54"'
$ws.Range("D46").Copy() | Out-Null
$ws.Range("D46").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("E46").Value = "  -0.17%  "

# Row 47
$ws.Range("E47").Value = "  -9.01%  "

# Row 48
$ws.Range("D48").Formula = '="0.0518"'
$ws.Range("D48").Copy() | Out-Null
$ws.Range("D48").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("E48").Value = "  +0.46%  "

# Row 49
$ws.Range("D49").Formula = '="19.01"'
$ws.Range("D49").Copy() | Out-Null
$ws.Range("D49").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("E49").Value = "  -6.93%  "

# Row 50
$ws.Range("D50").Formula = '="0.567"'
$ws.Range("D50").Copy() | Out-Null
$ws.Range("D50").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("E50").Value = "  -1.87%  "

# Row 51
$ws.Range("E51").Value = "  -2.87%  "
